$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.008931159973145
$ws.Range("B1").Value = 1.221996903419495
$ws.Range("C1").Value = 1.650306820869446
$ws.Range("D1").Value = 3.165316581726074
$ws.Range("E1").Value = 2.506096601486206
